$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F15").Value = "COSME0014SavePropsectRequest"
$ws.Range("E15").Value = "fullname = InvalidFields"
$ws.Range("G15").Value = "COSME0014SavePropsectResponse"
$ws.Range("C15").Value = "SME Prospect Invalid"
$ws.Range("D15").Value = "SME Prospect Invalid"

$ws.Range("C15").Select()
